$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.814.41"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.109.63"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "388.67"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").Value = "103.76"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "37.22"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D12").Value = "0.0862"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "3.598.68"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "18.75"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "7.91"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "3.109.72"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "10.71"
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").Value = "51.900.42"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("D21").Value = "12.52"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "'70.00"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "268.72"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").Value = "8.11"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").Value = "27.16"
$ws.Range("E27").Value = "  +3.81%  "
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "10.38"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "'35.50"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "50.39"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "3.42"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("D39").Value = "0.289"
$ws.Range("E39").Value = "  +6.49%  "
$ws.Range("D40").Value = "'1.90"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.60"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'16.90"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Value = "128.95"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("D46").Value = "22.34"
$ws.Range("E46").Value = "  +4.57%  "
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +6.86%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").Value = "2.048.96"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "3.416.77"
$ws.Range("D51").Value = "0.206"
$ws.Range("E51").Value = "  +5.47%  "
